# Adds a new "StringsCollectionWorksheet" worksheet that demonstrates a
# headerless / simple-type collection export (see commit message:
# "Added WithCollection method that assumes headerless, simple types").
#
# Layout:
#       A          B            C                   D
#   1  LastName   First Name   LanguagesSpoken1    LanguagesSpoken2
#   2  Test       Name         Spanish             Romanian
#   3  Test       Name2        English
#   4  Test       Name2

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new worksheet as the LAST tab in the workbook.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "StringsCollectionWorksheet"

# ---------------------------------------------------------------------
# 2. Write the grid values (plain strings, no numbers/dates). The write
#    order below mirrors how the cells were first authored (new text
#    values are interned in first-seen order).
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "LastName"
$ws.Cells.Item(1, 2).Value = "First Name"
$ws.Cells.Item(2, 1).Value = "Test"
$ws.Cells.Item(1, 3).Value = "LanguagesSpoken1"
$ws.Cells.Item(1, 4).Value = "LanguagesSpoken2"
$ws.Cells.Item(2, 3).Value = "Spanish"
$ws.Cells.Item(2, 4).Value = "Romanian"
$ws.Cells.Item(3, 2).Value = "Name2"
$ws.Cells.Item(3, 3).Value = "English"

$ws.Cells.Item(2, 2).Value = "Name"
$ws.Cells.Item(3, 1).Value = "Test"
$ws.Cells.Item(4, 1).Value = "Test"
$ws.Cells.Item(4, 2).Value = "Name2"

# ---------------------------------------------------------------------
# 3. Column widths / row height to match the generated report look.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 11.8776
$ws.Columns.Item(3).ColumnWidth = 8.1667
$ws.Columns.Item(4).ColumnWidth = 9.0221
$ws.Rows.Item(1).RowHeight = 25.5

# ---------------------------------------------------------------------
# 4. Header row styling: bold white text on a dark fill, thin light
#    border, left/top aligned with wrap, stored as text ("@").
# ---------------------------------------------------------------------
$headerRange = $ws.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.Font.Size = 10
$headerRange.Font.ThemeColor = 0
$headerRange.Interior.ThemeColor = 1
$headerRange.NumberFormat = "@"
$headerRange.HorizontalAlignment = -4131   # xlLeft
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.WrapText = $true
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2
$headerRange.Borders.Color = 14277081      # light grey (theme bg1, tinted)

# ---------------------------------------------------------------------
# 5. Body styling: normal weight font, left/top aligned with wrap.
#    The first data row additionally gets the thin table-style border.
# ---------------------------------------------------------------------
$bodyRange = $ws.Range("A2:D4")
$bodyRange.Font.Size = 11
$bodyRange.HorizontalAlignment = -4131     # xlLeft
$bodyRange.VerticalAlignment = -4160       # xlTop
$bodyRange.WrapText = $true

$firstDataRow = $ws.Range("A2:C2")
$firstDataRow.Borders.LineStyle = 1
$firstDataRow.Borders.Weight = 2
$firstDataRow.Borders.Color = 14277081

# ---------------------------------------------------------------------
# 6. Selection / activation so this becomes the active tab, matching
#    the last-used cell in the source workbook.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("C4").Select()
